$wb = $excel.ActiveWorkbook

# Both the "展览" sheet and the "全部类型" sheet received identical data updates.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Small "want to go" count corrections on existing rows ---
    $ws.Cells.Item(2, 6).Value = 6450   # F2: 6432 -> 6450
    $ws.Cells.Item(5, 6).Value = 40     # F5: 39 -> 40
    $ws.Cells.Item(6, 6).Value = 1931   # F6: 1924 -> 1931
    $ws.Cells.Item(7, 6).Value = 1467   # F7: 1465 -> 1467
    $ws.Cells.Item(9, 6).Value = 988    # F9: 987 -> 988
    $ws.Cells.Item(10, 6).Value = 323   # F10: 317 -> 323

    # --- Insert a new event as row 11, pushing the old rows 11-12 down to 12-13 ---
    # Column A only ever holds a plain running sequence number (0,1,2,...) that
    # is unrelated to which event occupies the row, so it is left untouched for
    # rows 11-12 and simply extended for the new row 13. Only columns B:I (the
    # actual event details) are shifted down. Copy (rather than read/write the
    # .Value) is used so literal text values (e.g. date-looking strings) are
    # preserved verbatim instead of being re-interpreted by Excel's smart data
    # typing.
    $ws.Cells.Item(12, 1).Copy($ws.Cells.Item(13, 1))
    $ws.Cells.Item(13, 1).Value = 12

    $ws.Range("B12:I12").Copy($ws.Range("B13:I13"))
    $ws.Range("B11:I11").Copy($ws.Range("B12:I12"))

    # The "want to go" count for the shifted "第二届漫画城市动漫展" event changed too.
    $ws.Cells.Item(12, 6).Value = 5606  # F12: 5603 -> 5606

    # --- Fill in the brand-new row 11 with the new event's details ---
    # Column B holds a literal date-like string ("2024-03-23"); prefix it with
    # an apostrophe so Excel stores it as text instead of auto-converting it
    # into a date serial number, then restore the plain "Normal" style so the
    # cell doesn't keep a lingering quote-prefix format marker.
    $ws.Cells.Item(11, 2).Value = "'2024-03-23"
    $ws.Cells.Item(11, 2).Style = "Normal"
    $ws.Cells.Item(11, 3).Value = "合肥·原神&星穹&崩铁only"
    $ws.Cells.Item(11, 4).Value = "金寨路与天堂窄路交叉口 梵木艺术中心"
    $ws.Cells.Item(11, 5).Value = "2024.03.23 09:00-03.23 17:00"
    $ws.Cells.Item(11, 6).Value = 1
    $ws.Cells.Item(11, 7).Value = 58
    $ws.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81574"
    $ws.Cells.Item(11, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/0V5uyX6C1706697212904.png"
}

$wb.Save()
